$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new blank column before column D (shifts D:K -> E:L).
$ws.Columns("D").Insert()

# 2. Give the new column D the same number formats / styles as column E
#    (the column it was inserted in front of) so per-row styles (date /
#    thousands number formats) line up with the rest of the table.
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# 3. Match the column width Excel used for the new column.
$ws.Columns("D").ColumnWidth = $ws.Columns("E").ColumnWidth

# 4. Populate the new column D with the refreshed period's figures.
#    Rows not listed here are the blank/ separator rows, where column D
#    is simply left empty (matching the source diff).
$newColumnD = @{
    7 = 43465;
    8 = 2947100;
    9 = 1631700;
    10 = 1315400;
    12 = 185800;
    13 = 0;
    14 = 907500;
    15 = 0;
    17 = 3371100;
    18 = -424000;
    20 = 20000;
    21 = 319700;
    22 = 534900;
    23 = -938800;
    24 = 17300;
    25 = 0;
    26 = -956100;
    27 = -956100;
    28 = 0;
    29 = -75400;
    30 = 0;
    31 = 0;
    32 = -20000;
    33 = -1031500;
    34 = 0;
    35 = -1031500;
    38 = 43465;
    41 = 1149100;
    43 = 510400;
    44 = 322200;
    45 = 361500;
    46 = 2343200;
    47 = 700;
    48 = 498900;
    49 = 7221900;
    50 = 0;
    51 = 0;
    52 = 67700;
    53 = 0;
    54 = 10132400;
    57 = 96000;
    58 = 34100;
    59 = 1819900;
    60 = 1950100;
    61 = 8224300;
    62 = 456300;
    63 = 0;
    64 = 0;
    65 = 0;
    66 = 10630700;
    68 = 0;
    69 = 0;
    70 = 0;
    71 = 0;
    72 = -9124900;
    73 = 0;
    74 = 0;
    75 = 0;
    76 = -498300;
    77 = 0;
    80 = 43465;
    81 = -1031500;
    83 = 723700;
    84 = 0;
    85 = 0;
    86 = 0;
    87 = 0;
    88 = 0;
    89 = 267300;
    92 = 0;
    93 = 0;
    94 = -17900;
    96 = 0;
    97 = 0;
    98 = 0;
    99 = 0;
    100 = -81600;
    101 = -2000;
    102 = 165800;
}

foreach ($row in $newColumnD.Keys) {
    $ws.Cells.Item([int]$row, 4).Value = $newColumnD[$row]
}

# Row 42 ("Short Term Investments") uses the literal "NA" text marker,
# matching the rest of that row.
$ws.Cells.Item(42, 4).Value = "NA"

# 5. Row 91 ("Capital Expenditures") was re-keyed with a new set of
#    figures across the whole row (not just a plain shift), so overwrite
#    D91:K91 explicitly.
$ws.Cells.Item(91, 4).Value = -83400
$ws.Cells.Item(91, 5).Value = -125700
$ws.Cells.Item(91, 6).Value = -138900
$ws.Cells.Item(91, 7).Value = -81800
$ws.Cells.Item(91, 8).Value = -80400
$ws.Cells.Item(91, 9).Value = -96500
$ws.Cells.Item(91, 10).Value = -99800
$ws.Cells.Item(91, 11).Value = -61700
